# Sample Project / Main.xlsx - "Rules" sheet update.
# Row 11 / column B ("R40") becomes the text "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

# Force the new value to be stored as text (matching the shared-string
# cell type used by the rest of the column) instead of being
# auto-detected as a number.
$cell.NumberFormat = "@"
$cell.Value = "1"
